$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Instructions")
$ws.Range("B3").Value = "ihcc-browser@googlegroups.com"
